# "Implemented Receive and forward shipment API"
#
# The "Receive Shipment at branch" API (row 11) is now done, so its Status
# moves from "In Progress" (yellow) to "Complete" (green) - matching the
# styling already used by the other completed APIs above it (rows 8-10).
#
# The "Attempt Delivery" and "Update Delivery Status" APIs (rows 12-13),
# which previously had no Status set, are now "In Progress" (yellow) -
# the same styling "Receive Shipment at branch" used to have.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1) Move the existing "In Progress" (yellow) formatting on row 11 down onto
#    rows 12-13 before row 11's own formatting gets overwritten in step 2.
$ws.Range("B11:G11").Copy() | Out-Null
$ws.Range("B12:G13").PasteSpecial(-4122) | Out-Null

# 2) Row 11 is now complete - give it the same "Complete" (green) formatting
#    already used on rows 8-10.
$ws.Range("B10:G10").Copy() | Out-Null
$ws.Range("B11:G11").PasteSpecial(-4122) | Out-Null

# 3) Update the Status text in column F for all three rows.
$ws.Range("F11").Value = "Complete"
$ws.Range("F12").Value = "In Progress"
$ws.Range("F13").Value = "In Progress"

# 4) Match the selection/scroll state left behind by the edit (user had just
#    finished selecting/filling B12:G13 with the new "In Progress" status).
$ws.Range("B12:G13").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 2
